# Appends 16 new lesson rows (119-134) to the media-data tracker sheet:
#   "JavaScript - Operators" (3), "JavaScript - Object" (7) and
#   "JavaScript - External Script" (6) topics, per the commit message
#   "included the operators, object, external script".
#
# Columns are id(A) / topic(B) / seq(C) / subTopic(D) / media(E). The writes
# below are sequenced deliberately (not simply row-by-row) so that the
# resulting xl/sharedStrings.xml table gets new unique strings appended in
# the same order the original authored workbook has them in.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: row ids (plain numbers) ---
$ws.Range("A120").Value = 119
$ws.Range("A121").Value = 120
$ws.Range("A122").Value = 121
$ws.Range("A123").Value = 122
$ws.Range("A124").Value = 123
$ws.Range("A125").Value = 124
$ws.Range("A126").Value = 125
$ws.Range("A127").Value = 126
$ws.Range("A128").Value = 127
$ws.Range("A129").Value = 128
$ws.Range("A130").Value = 129
$ws.Range("A131").Value = 130
$ws.Range("A132").Value = 131
$ws.Range("A133").Value = 132
$ws.Range("A134").Value = 133
$ws.Range("A135").Value = 134

# --- Unique-value cells, in shared-string insertion order ---
# (seq values in column C are entered with a leading apostrophe so they are
# stored as text with the existing quote-prefix style, matching the sheet's
# other seq cells, instead of being auto-converted to numbers.)
$ws.Range("B120").Value = 'JavaScript - Operators'
$ws.Range("C120").Value = "'1501"
$ws.Range("D120").Value = 'Assignment Operator'
$ws.Range("D121").Value = 'Arithmetic Operator'
$ws.Range("D122").Value = 'Comparison Operator'
$ws.Range("C121").Value = "'1502"
$ws.Range("C122").Value = "'1503"
$ws.Range("B123").Value = 'JavaScript - Object'
$ws.Range("C123").Value = "'1601"
$ws.Range("D123").Value = 'Limitations of Variables and Arrays'
$ws.Range("D124").Value = 'Display all menu item details in console'
$ws.Range("D126").Value = 'Implement Object Array'
$ws.Range("D125").Value = 'Implement Object'
$ws.Range("D127").Value = 'Menu Object Array Walkthrough'
$ws.Range("D128").Value = 'Include Menu Object Array'
$ws.Range("D129").Value = 'Display menu based on object array'
$ws.Range("C124").Value = "'1602"
$ws.Range("C125").Value = "'1603"
$ws.Range("C126").Value = "'1604"
$ws.Range("C127").Value = "'1605"
$ws.Range("C128").Value = "'1606"
$ws.Range("C129").Value = "'1607"
$ws.Range("B130").Value = 'JavaScript - External Script'
$ws.Range("C130").Value = "'1701"
$ws.Range("D130").Value = 'External JavaScript'
$ws.Range("D131").Value = 'Copy menu object array script file'
$ws.Range("D132").Value = 'Create HTML for learning external script'
$ws.Range("D133").Value = 'Remove array declaration'
$ws.Range("D134").Value = 'Include script tag'
$ws.Range("D135").Value = 'Script tags code execution flow'
$ws.Range("C131").Value = "'1702"
$ws.Range("C132").Value = "'1703"
$ws.Range("C133").Value = "'1704"
$ws.Range("C134").Value = "'1705"
$ws.Range("C135").Value = "'1706"
$ws.Range("E120").Value = 'https://www.youtube.com/embed/fcRHTqo9WvU'
$ws.Range("E121").Value = 'https://www.youtube.com/embed/JqtHkK0Tc8c'
$ws.Range("E122").Value = 'https://www.youtube.com/embed/YRkleNBoGi0'
$ws.Range("E123").Value = 'https://www.youtube.com/embed/CjGEW-y_F5k'
$ws.Range("E124").Value = 'https://www.youtube.com/embed/-HRBaPd9l5Q'
$ws.Range("E125").Value = 'https://www.youtube.com/embed/VtMOkEAN5II'
$ws.Range("E126").Value = 'https://www.youtube.com/embed/oZ_TcTAUDco'
$ws.Range("E127").Value = 'https://www.youtube.com/embed/QIMZ8TrEIXg'
$ws.Range("E128").Value = 'https://www.youtube.com/embed/FOMuVv3mWXs'
$ws.Range("E129").Value = 'https://www.youtube.com/embed/e79er5dw_Eo'
$ws.Range("E130").Value = 'https://www.youtube.com/embed/zNfz68dFlRE'
$ws.Range("E131").Value = 'https://www.youtube.com/embed/CiYlUmaTsjs'
$ws.Range("E132").Value = 'https://www.youtube.com/embed/i-hUchsNmCw'
$ws.Range("E133").Value = 'https://www.youtube.com/embed/VM9yyE3_WT0'
$ws.Range("E134").Value = 'https://www.youtube.com/embed/JtMkzEh4q6E'
$ws.Range("E135").Value = 'https://www.youtube.com/embed/ZLW4FxYBX28'

# --- Remaining repeated Topic (column B) cells ---
$ws.Range("B121").Value = 'JavaScript - Operators'
$ws.Range("B122").Value = 'JavaScript - Operators'
$ws.Range("B124").Value = 'JavaScript - Object'
$ws.Range("B125").Value = 'JavaScript - Object'
$ws.Range("B126").Value = 'JavaScript - Object'
$ws.Range("B127").Value = 'JavaScript - Object'
$ws.Range("B128").Value = 'JavaScript - Object'
$ws.Range("B129").Value = 'JavaScript - Object'
$ws.Range("B131").Value = 'JavaScript - External Script'
$ws.Range("B132").Value = 'JavaScript - External Script'
$ws.Range("B133").Value = 'JavaScript - External Script'
$ws.Range("B134").Value = 'JavaScript - External Script'
$ws.Range("B135").Value = 'JavaScript - External Script'

# --- Final selection matches the authored view state ---
$ws.Range("A119").Select()
